$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 490, shifting the
# existing rows 490-524 down to 492-526 (dimension grows from R524 to R526).
$ws.Rows.Item(490).Insert()
$ws.Rows.Item(491).Insert()

# New row 490: Papa / Patagonia / 1a (cosecha) record for 2022-07-04
$ws.Cells.Item(490, 1).Value = 5
$ws.Cells.Item(490, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(490, 3).Value = "Maule"
$ws.Cells.Item(490, 4).Value = 44746
$ws.Cells.Item(490, 5).Value = 7
$ws.Cells.Item(490, 6).Value = 100114001
$ws.Cells.Item(490, 7).Value = "Papa"
$ws.Cells.Item(490, 8).Value = "Patagonia"
$ws.Cells.Item(490, 9).Value = "1a (cosecha)"
$ws.Cells.Item(490, 10).Value = 1200
$ws.Cells.Item(490, 11).Value = 7000
$ws.Cells.Item(490, 12).Value = 7000
$ws.Cells.Item(490, 13).Value = 7000
$ws.Cells.Item(490, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(490, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(490, 16).Value = 280
$ws.Cells.Item(490, 17).Value = 25
$ws.Cells.Item(490, 18).Value = "Hortaliza"

# New row 491: Papa / Yagana / 1a (cosecha) record for 2022-07-04
$ws.Cells.Item(491, 1).Value = 5
$ws.Cells.Item(491, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(491, 3).Value = "Maule"
$ws.Cells.Item(491, 4).Value = 44746
$ws.Cells.Item(491, 5).Value = 7
$ws.Cells.Item(491, 6).Value = 100114001
$ws.Cells.Item(491, 7).Value = "Papa"
$ws.Cells.Item(491, 8).Value = "Yagana"
$ws.Cells.Item(491, 9).Value = "1a (cosecha)"
$ws.Cells.Item(491, 10).Value = 500
$ws.Cells.Item(491, 11).Value = 7500
$ws.Cells.Item(491, 12).Value = 7500
$ws.Cells.Item(491, 13).Value = 7500
$ws.Cells.Item(491, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(491, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(491, 16).Value = 300
$ws.Cells.Item(491, 17).Value = 25
$ws.Cells.Item(491, 18).Value = "Hortaliza"
